$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

# --- Add new row 23 -------------------------------------------------------
# Copy the formatting of row 22 (columns B:G, style index 9) down onto row 23
# before writing values, so the new row picks up the same cell style that
# Excel would have auto-extended from the table/autofilter region.
$ws.Range("B22:G22").Copy()
$ws.Range("B23:G23").PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item(23, 1).Value = 23
$ws.Cells.Item(23, 2).Value = "Create a detailed syllabus for communication"
$ws.Cells.Item(23, 3).Value = "Marketing"
$ws.Cells.Item(23, 4).Value = "Rahul"
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(23, 6).Value = "WIP"

# --- Re-apply the AutoFilter over the grown range, filtered to Status=WIP -
# Clear the existing autofilter first so the new range (A1:H23) "sticks"
# instead of Excel keeping the old A1:H22 extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:H23").AutoFilter(6, @("WIP"), 7)

# --- Keep the workbook-level _FilterDatabase defined name in sync ---------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Tracker!_FilterDatabase") {
        $n.RefersTo = "=Tracker!`$A`$1:`$H`$23"
    }
}

# --- Restore the active selection to B16 ----------------------------------
$ws.Activate()
$ws.Range("B16").Select()
